# Identifiants.xlsx -- "Ajout des icones et fin du programme Sigfox"
#
# 1) App EUI / App Key hash strings (B3/B4) and the Device-EUI list (B13:B19)
#    become =UPPER("...") formulas instead of literal lower-case text.
# 2) The little "card" summary (rows 9-10, column B) gets Text number format.
# 3) The "Arduino" table's second column is renamed "Device EUI" -> "Device EUI 16byte".
# 4) A new copy of the App EUI / App Key card is appended at rows 22-24,
#    holding the short (16 hex-digit) uppercase forms of the two hashes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoRa")

# --- 1) Upper-case the long identifiers via formulas -----------------------
$ws.Range("B3").Formula = '=UPPER("50f925fb23b68de3b85a1fcf55c989de7d6a57c1")'
$ws.Range("B4").Formula = '=UPPER("6475463840e57f1dc21360a72f14d2a1c25b2406")'

$ws.Range("B13").Formula = '=UPPER("a8610a30393d6d05")'
$ws.Range("B14").Formula = '=UPPER("a8610a30393a6605")'
$ws.Range("B15").Formula = '=UPPER("a8610a3039316905")'
$ws.Range("B16").Formula = '=UPPER("a8610a3039246c05")'
$ws.Range("B17").Formula = '=UPPER("a8610a30393e7205")'
$ws.Range("B18").Formula = '=UPPER("a8610a30393f7705")'
$ws.Range("B19").Formula = '=UPPER("a8610a30392f7805")'

# --- 2) Card at rows 8-10: column B now stores its value as Text -----------
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B10").NumberFormat = "@"

# --- 3) Rename the Arduino/Device-EUI table header -------------------------
$ws.Range("B12").Value = "Device EUI 16byte"

# --- 4) New card at rows 22-24 (App EUI / App Key, 16-byte uppercase) ------
$ws.Range("A8:B8").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A9:B9").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A10:B10").Copy()
$ws.Range("A24").PasteSpecial(-4122)

$ws.Range("A22").Value = "Type d'identifiant"
$ws.Range("B22").Value = "Identifiant"
$ws.Range("A23").Value = "App EUI"
$ws.Range("B23").Value = "50F925FB23B68DE3"
$ws.Range("A24").Value = "App Key"
$ws.Range("B24").Value = "6475463840E57F1DC21360A72F14D2A1"

# --- View bits: zoom + selection -------------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("C3").Select()
